$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 46074 -> 46075) for every
# data row (rows 2 through 239). Increment each by one day.
for ($r = 2; $r -le 239; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
